$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Atualizei dados da bibi (faturamento anual - linha 2025 / row 9)
$ws.Range("B9").Value = 4349777.58
$ws.Range("C9").Value = 696076.24
$ws.Range("D9").Value = 5045853.82
$ws.Range("E9").Value = 13.79501398239079
$ws.Range("F9").Value = 86.2049860176092
$ws.Range("G9").Value = -32.7277022001421
$ws.Range("H9").Value = -21.44902934008548
$ws.Range("I9").Value = 43517
$ws.Range("J9").Value = 1889
$ws.Range("K9").Value = 45406
$ws.Range("L9").Value = 31583
$ws.Range("M9").Value = 159.7648678086312
$ws.Range("N9").Value = 9.074729300549267
